$wb = $excel.ActiveWorkbook

# --- Add "imports" worksheet after the last existing sheet ("html") ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$imports = $wb.Worksheets.Add($null, $last)
$imports.Name = "imports"
$imports.Range("A1").Value = "model_name"
$imports.Range("B1").Value = "functions"
$imports.Range("A1:B1").Font.Color = 0
[void]$imports.Range("A1:B1").Select()

# --- Add "admin" worksheet after "imports" ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$admin = $wb.Worksheets.Add($null, $last)
$admin.Name = "admin"
$admin.Range("A1").Value = "model_name"
$admin.Range("B1").Value = "functions"
$admin.Range("A1:B1").Font.Color = 0
[void]$admin.Range("A1:B1").Select()

# --- Add "apps" worksheet after "admin" ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$apps = $wb.Worksheets.Add($null, $last)
$apps.Name = "apps"
$apps.Range("A1").Value = "model_name"
$apps.Range("B1").Value = "functions"
$apps.Range("A1:B1").Font.Color = 0
[void]$apps.Range("B9").Select()

$apps.Activate()
